$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1856
$ws.Range("I92").Value = 1677.2
$ws.Range("K92").Value = 1677.2
$ws.Range("M92").Value = -429.2
$ws.Range("H96").Value = 1384.375
$ws.Range("I96").Value = 3542
$ws.Range("J96").Value = 886.46155
$ws.Range("K96").Value = 10626
$ws.Range("L96").Value = 2659.38465
$ws.Range("M96").Value = -9253
$ws.Range("N96").Value = -5405.38465
$ws.Range("H99").Value = 101682.5
$ws.Range("I99").Value = 790.6667
$ws.Range("J99").Value = 253020.25
$ws.Range("K99").Value = 2372.0001
$ws.Range("L99").Value = 759060.75
$ws.Range("M99").Value = -874.0001000000002
$ws.Range("N99").Value = -762056.75
$ws.Range("H138").Value = 3045
$ws.Range("I138").Value = 2289.3635
$ws.Range("J138").Value = 3800.6365
$ws.Range("K138").Value = 6868.0905
$ws.Range("L138").Value = 11401.9095
$ws.Range("M138").Value = -1728.0905
$ws.Range("N138").Value = -21681.9095
$ws.Range("H140").Value = 81240
$ws.Range("J140").Value = 81240
$ws.Range("L140").Value = 81240
$ws.Range("N140").Value = -91600
$ws.Range("H141").Value = 3591.8462
$ws.Range("I141").Value = 3417.182
$ws.Range("K141").Value = 10251.546
$ws.Range("M141").Value = -5071.545999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13890798
$ws.Range("I61").Value = 15153418
$ws.Range("J61").Value = 1971.3334
$ws.Range("K61").Value = 15153418
$ws.Range("L61").Value = 1971.3334
$ws.Range("M61").Value = -15153206
$ws.Range("N61").Value = -2395.3334
$ws.Range("H122").Value = 5634.84
$ws.Range("I122").Value = 6293.4287
$ws.Range("J122").Value = 2177.25
$ws.Range("K122").Value = 18880.2861
$ws.Range("L122").Value = 6531.75
$ws.Range("M122").Value = -16430.2861
$ws.Range("N122").Value = -11431.75
$ws.Range("H136").Value = 13890798
$ws.Range("I136").Value = 15153418
$ws.Range("J136").Value = 1971.3334
$ws.Range("K136").Value = 45460254
$ws.Range("L136").Value = 5914.0002
$ws.Range("M136").Value = -45457704
$ws.Range("N136").Value = -11014.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 411
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 444
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 444
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -790
$ws.Range("H94").Value = 984.5263
$ws.Range("I94").Value = 761
$ws.Range("J94").Value = 1414.3846
$ws.Range("K94").Value = 761
$ws.Range("L94").Value = 1414.3846
$ws.Range("M94").Value = -310
$ws.Range("N94").Value = -2316.3846
$ws.Range("H99").Value = 806.2941
$ws.Range("I99").Value = 762.9375
$ws.Range("K99").Value = 762.9375
$ws.Range("M99").Value = 735.0625
$ws.Range("H105").Value = 4197.2705
$ws.Range("I105").Value = 3183.5454
$ws.Range("J105").Value = 4626.154
$ws.Range("K105").Value = 3183.5454
$ws.Range("L105").Value = 4626.154
$ws.Range("M105").Value = -1436.5454
$ws.Range("N105").Value = -8120.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2660.125
$ws.Range("I58").Value = 596.63635
$ws.Range("J58").Value = 7199.8
$ws.Range("K58").Value = 596.63635
$ws.Range("L58").Value = 7199.8
$ws.Range("M58").Value = -393.63635
$ws.Range("N58").Value = -7605.8
$ws.Range("H132").Value = 2829.2593
$ws.Range("I132").Value = 1918.9
$ws.Range("J132").Value = 5430.2856
$ws.Range("K132").Value = 5756.700000000001
$ws.Range("L132").Value = 16290.8568
$ws.Range("M132").Value = -3226.700000000001
$ws.Range("N132").Value = -21350.8568
$ws.Range("H136").Value = 2660.125
$ws.Range("I136").Value = 596.63635
$ws.Range("J136").Value = 7199.8
$ws.Range("K136").Value = 1789.90905
$ws.Range("L136").Value = 21599.4
$ws.Range("M136").Value = 760.09095
$ws.Range("N136").Value = -26699.4
$ws.Range("H140").Value = 32177.777
$ws.Range("J140").Value = 32177.777
$ws.Range("L140").Value = 32177.777
$ws.Range("N140").Value = -42537.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 996.8333
$ws.Range("J17").Value = 996.4
$ws.Range("L17").Value = 2989.2
$ws.Range("N17").Value = -3327.2
$ws.Range("H80").Value = 2720
$ws.Range("J80").Value = 2775
$ws.Range("L80").Value = 8325
$ws.Range("N80").Value = -10197
$ws.Range("H83").Value = 2720
$ws.Range("J83").Value = 2775
$ws.Range("L83").Value = 24975
$ws.Range("N83").Value = -34335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1157.75
$ws.Range("J97").Value = 1308.25
$ws.Range("L97").Value = 1308.25
$ws.Range("N97").Value = -2300.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1667.2222
$ws.Range("I68").Value = 2110.5715
$ws.Range("J68").Value = 1385.091
$ws.Range("K68").Value = 2110.5715
$ws.Range("L68").Value = 1385.091
$ws.Range("M68").Value = -1361.5715
$ws.Range("N68").Value = -2883.091
$ws.Range("H71").Value = 1667.2222
$ws.Range("I71").Value = 2110.5715
$ws.Range("J71").Value = 1385.091
$ws.Range("K71").Value = 10552.8575
$ws.Range("L71").Value = 6925.455
$ws.Range("M71").Value = -6808.8575
$ws.Range("N71").Value = -14413.455
$ws.Range("H93").Value = 1322
$ws.Range("I93").Value = 1233.3077
$ws.Range("J93").Value = 1466.125
$ws.Range("K93").Value = 1233.3077
$ws.Range("L93").Value = 1466.125
$ws.Range("M93").Value = 14.69229999999993
$ws.Range("N93").Value = -3962.125
$ws.Range("H122").Value = 8300
$ws.Range("I122").Value = 10983.333
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 32949.999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -30499.999
$ws.Range("N122").Value = -22900
$ws.Range("H139").Value = 44481.273
$ws.Range("J139").Value = 44864.4
$ws.Range("L139").Value = 44864.4
$ws.Range("N139").Value = -55144.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1172.725
$ws.Range("I126").Value = 1326.4839
$ws.Range("J126").Value = 643.1111
$ws.Range("K126").Value = 3979.4517
$ws.Range("L126").Value = 1929.3333
$ws.Range("M126").Value = -1509.4517
$ws.Range("N126").Value = -6869.3333
